$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "M1"
$ws.Cells.Item(2,2).Value = "Spn"
$ws.Cells.Item(2,3).Value = "Siglec1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 3.787414
$ws.Cells.Item(2,8).Value = 11.362242
$ws.Cells.Item(2,9).Value = 0.5877125485801681
$ws.Cells.Item(2,10).Value = 0.587712548580168
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.120277
$ws.Cells.Item(2,14).Value = 0.360831
$ws.Cells.Item(2,15).Value = 0.001062914843064901
$ws.Cells.Item(2,16).Value = 0.00106414053856905
$ws.Cells.Item(2,17).Value = 0.455538793678
$ws.Cells.Item(2,18).Value = 4.099849143102
$ws.Cells.Item(2,19).Value = 0.0006246883913413623
$ws.Cells.Item(2,20).Value = 0.0006254087479698888

# Row 3
$ws.Cells.Item(3,1).Value = "M1"
$ws.Cells.Item(3,2).Value = "Spn"
$ws.Cells.Item(3,3).Value = "Siglec1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 3.787414
$ws.Cells.Item(3,8).Value = 11.362242
$ws.Cells.Item(3,9).Value = 0.5877125485801681
$ws.Cells.Item(3,10).Value = 0.587712548580168
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.4782236666666666
$ws.Cells.Item(3,14).Value = 1.434671
$ws.Cells.Item(3,15).Value = 0.004226169871254865
$ws.Cells.Item(3,16).Value = 0.004231043260167216
$ws.Cells.Item(3,17).Value = 1.811231010264666
$ws.Cells.Item(3,18).Value = 16.301079092382
$ws.Cells.Item(3,19).Value = 0.002483773065767917
$ws.Cells.Item(3,20).Value = 0.002486637217585817

# Row 4
$ws.Cells.Item(4,1).Value = "M1"
$ws.Cells.Item(4,2).Value = "Spn"
$ws.Cells.Item(4,3).Value = "Siglec1"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 3.787414
$ws.Cells.Item(4,8).Value = 11.362242
$ws.Cells.Item(4,9).Value = 0.5877125485801681
$ws.Cells.Item(4,10).Value = 0.587712548580168
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 69.61810866666667
$ws.Cells.Item(4,14).Value = 208.854326
$ws.Cells.Item(4,15).Value = 0.6152308508518273
$ws.Cells.Item(4,16).Value = 0.615940301559777
$ws.Cells.Item(4,17).Value = 263.6725994176547
$ws.Cells.Item(4,18).Value = 2373.053394758892
$ws.Cells.Item(4,19).Value = 0.3615788913192727
$ws.Cells.Item(4,20).Value = 0.3619958444029338

# Row 5
$ws.Cells.Item(5,1).Value = "M1"
$ws.Cells.Item(5,2).Value = "Spn"
$ws.Cells.Item(5,3).Value = "Siglec1"
$ws.Cells.Item(5,4).Value = "M2"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 3.787414
$ws.Cells.Item(5,8).Value = 11.362242
$ws.Cells.Item(5,9).Value = 0.5877125485801681
$ws.Cells.Item(5,10).Value = 0.587712548580168
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 42.55008066666667
$ws.Cells.Item(5,14).Value = 127.650242
$ws.Cells.Item(5,15).Value = 0.376024612471286
$ws.Cells.Item(5,16).Value = 0.376458223573777
$ws.Cells.Item(5,17).Value = 161.1547712180627
$ws.Cells.Item(5,18).Value = 1450.392940962564
$ws.Cells.Item(5,19).Value = 0.2209943833243696
$ws.Cells.Item(5,20).Value = 0.2212492220105071

# Row 6
$ws.Cells.Item(6,1).Value = "M1"
$ws.Cells.Item(6,2).Value = "Spn"
$ws.Cells.Item(6,3).Value = "Siglec1"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 3.787414
$ws.Cells.Item(6,8).Value = 11.362242
$ws.Cells.Item(6,9).Value = 0.5877125485801681
$ws.Cells.Item(6,10).Value = 0.587712548580168
$ws.Cells.Item(6,11).Value = 2
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 0.391011
$ws.Cells.Item(6,14).Value = 0.782022
$ws.Cells.Item(6,15).Value = 0.003455451962566825
$ws.Cells.Item(6,16).Value = 0.00230629106770994
$ws.Cells.Item(6,17).Value = 1.480920535554
$ws.Cells.Item(6,18).Value = 8.885523213323999
$ws.Cells.Item(6,19).Value = 0.002030812479416492
$ws.Cells.Item(6,20).Value = 0.001355436201171486

# Row 7
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Spn"
$ws.Cells.Item(7,3).Value = "Siglec1"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 2.656916666666667
$ws.Cells.Item(7,8).Value = 7.97075
$ws.Cells.Item(7,9).Value = 0.412287451419832
$ws.Cells.Item(7,10).Value = 0.4122874514198319
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.120277
$ws.Cells.Item(7,14).Value = 0.360831
$ws.Cells.Item(7,15).Value = 0.001062914843064901
$ws.Cells.Item(7,16).Value = 0.00106414053856905
$ws.Cells.Item(7,17).Value = 0.3195659659166667
$ws.Cells.Item(7,18).Value = 2.87609369325
$ws.Cells.Item(7,19).Value = 0.0004382264517235387
$ws.Cells.Item(7,20).Value = 0.0004387317905991608

# Row 8
$ws.Cells.Item(8,1).Value = "M2"
$ws.Cells.Item(8,2).Value = "Spn"
$ws.Cells.Item(8,3).Value = "Siglec1"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 2.656916666666667
$ws.Cells.Item(8,8).Value = 7.97075
$ws.Cells.Item(8,9).Value = 0.412287451419832
$ws.Cells.Item(8,10).Value = 0.4122874514198319
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.4782236666666666
$ws.Cells.Item(8,14).Value = 1.434671
$ws.Cells.Item(8,15).Value = 0.004226169871254865
$ws.Cells.Item(8,16).Value = 0.004231043260167216
$ws.Cells.Item(8,17).Value = 1.270600430361111
$ws.Cells.Item(8,18).Value = 11.43540387325
$ws.Cells.Item(8,19).Value = 0.001742396805486948
$ws.Cells.Item(8,20).Value = 0.001744406042581398

# Row 9
$ws.Cells.Item(9,1).Value = "M2"
$ws.Cells.Item(9,2).Value = "Spn"
$ws.Cells.Item(9,3).Value = "Siglec1"
$ws.Cells.Item(9,4).Value = "M1"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 2.656916666666667
$ws.Cells.Item(9,8).Value = 7.97075
$ws.Cells.Item(9,9).Value = 0.412287451419832
$ws.Cells.Item(9,10).Value = 0.4122874514198319
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 69.61810866666667
$ws.Cells.Item(9,14).Value = 208.854326
$ws.Cells.Item(9,15).Value = 0.6152308508518273
$ws.Cells.Item(9,16).Value = 0.615940301559777
$ws.Cells.Item(9,17).Value = 184.9695132182778
$ws.Cells.Item(9,18).Value = 1664.7256189645
$ws.Cells.Item(9,19).Value = 0.2536519595325546
$ws.Cells.Item(9,20).Value = 0.2539444571568432

# Row 10
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Spn"
$ws.Cells.Item(10,3).Value = "Siglec1"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 2.656916666666667
$ws.Cells.Item(10,8).Value = 7.97075
$ws.Cells.Item(10,9).Value = 0.412287451419832
$ws.Cells.Item(10,10).Value = 0.4122874514198319
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 42.55008066666667
$ws.Cells.Item(10,14).Value = 127.650242
$ws.Cells.Item(10,15).Value = 0.376024612471286
$ws.Cells.Item(10,16).Value = 0.376458223573777
$ws.Cells.Item(10,17).Value = 113.0520184912778
$ws.Cells.Item(10,18).Value = 1017.4681664215
$ws.Cells.Item(10,19).Value = 0.1550302291469165
$ws.Cells.Item(10,20).Value = 0.1552090015632698

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Spn"
$ws.Cells.Item(11,3).Value = "Siglec1"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 2.656916666666667
$ws.Cells.Item(11,8).Value = 7.97075
$ws.Cells.Item(11,9).Value = 0.412287451419832
$ws.Cells.Item(11,10).Value = 0.4122874514198319
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 0.391011
$ws.Cells.Item(11,14).Value = 0.782022
$ws.Cells.Item(11,15).Value = 0.003455451962566825
$ws.Cells.Item(11,16).Value = 0.00230629106770994
$ws.Cells.Item(11,17).Value = 1.03888364275
$ws.Cells.Item(11,18).Value = 6.2333018565
$ws.Cells.Item(11,19).Value = 0.001424639483150333
$ws.Cells.Item(11,20).Value = 0.0009508548665384542
